$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'28.623.98"
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -2.13%  '
$ws.Range('D3').Value = "'1.792.73"
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -2.01%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'231.72"
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.66%  '
$ws.Range('D6').Value = "'0.5880"
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.42%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = "'0.2766"
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.94%  '
$ws.Range('D9').Value = "'0.06739"
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -4.24%  '
$ws.Range('D10').Value = "'23.15"
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -1.72%  '
$ws.Range('E11').Value = '  -1.59%  '
$ws.Range('D12').Value = "'1.788.69"
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -2.11%  '
$ws.Range('D13').Value = "'4.778"
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -0.21%  '
$ws.Range('D14').Value = "'0.6122"
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -2.46%  '
$ws.Range('D15').Value = "'2.036.35"
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -1.98%  '
$ws.Range('D16').Value = "'75.26"
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -4.70%  '
$ws.Range('D17').Value = "'0.000008783"
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -10.06%  '
$ws.Range('D18').Value = "'28.611.39"
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -2.23%  '
$ws.Range('D19').Value = "'5.403"
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -7.12%  '
$ws.Range('E20').Value = '  -0.01%  '
$ws.Range('D21').Value = "'208.33"
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -7.04%  '
$ws.Range('E22').Value = '  -1.87%  '
$ws.Range('D23').Value = "'6.822"
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -2.51%  '
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('D25').Value = "'152.93"
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -2.12%  '
$ws.Range('D26').Value = "'8.130"
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +1.85%  '
$ws.Range('D27').Value = "'0.1256"
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -3.40%  '
$ws.Range('D28').Value = "'16.38"
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -1.36%  '
$ws.Range('D29').Value = "'1.410"
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -4.02%  '
$ws.Range('D30').Value = "'0.06222"
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -6.33%  '
$ws.Range('E31').Value = '  -1.74%  '
$ws.Range('D32').Value = "'3.802"
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.36%  '
$ws.Range('E33').Value = '  -1.63%  '
$ws.Range('D34').Value = "'1.738"
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.91%  '
$ws.Range('E35').Value = '  -5.28%  '
$ws.Range('D36').Value = "'0.6384"
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -1.08%  '
$ws.Range('D37').Value = "'2.500"
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -1.81%  '
$ws.Range('D38').Value = "'2.714"
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.84%  '
$ws.Range('D39').Value = "'0.01696"
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -3.43%  '
$ws.Range('D40').Value = "'6.365"
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -2.42%  '
$ws.Range('D41').Value = "'1.141.29"
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -5.81%  '
$ws.Range('D42').Value = "'0.8742"
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -3.05%  '
$ws.Range('E43').Value = '  -0.11%  '
$ws.Range('D44').Value = "'100.29"
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.23%  '
$ws.Range('D45').Value = "'1.946.03"
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -2.19%  '
$ws.Range('D46').Value = "'59.88"
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -4.26%  '
$ws.Range('D47').Value = "'0.00000000112"
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -3.40%  '
$ws.Range('D48').Value = "'1.582"
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.05%  '
$ws.Range('D49').Value = "'8.344"
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -2.37%  '
$ws.Range('D50').Value = "'0.05468"
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.60%  '
$ws.Range('D51').Value = "'0.4487"
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.42%  '
